# Apply diff: update rows 442-526 (shift pattern) and add new rows 527-529
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 442
$ws.Range("D442").Value = 44785
$ws.Range("L442").Value = "Especial"
$ws.Range("M442").Value = 200
$ws.Range("N442").Value = 35000
$ws.Range("O442").Value = 35000
$ws.Range("P442").Value = 35000
$ws.Range("Q442").Value = "`$/bandeja 10 kilos"
$ws.Range("R442").Value = "Perú"
$ws.Range("S442").Value = 3500
$ws.Range("T442").Value = 10

# Row 443
$ws.Range("D443").Value = 44785
$ws.Range("L443").Value = "Primera"
$ws.Range("M443").Value = 200
$ws.Range("N443").Value = 33000
$ws.Range("O443").Value = 33000
$ws.Range("P443").Value = 33000
$ws.Range("Q443").Value = "`$/bandeja 10 kilos"
$ws.Range("R443").Value = "Perú"
$ws.Range("S443").Value = 3300
$ws.Range("T443").Value = 10

# Row 444
$ws.Range("D444").Value = 44785
$ws.Range("L444").Value = "Segunda"
$ws.Range("M444").Value = 200
$ws.Range("N444").Value = 30000
$ws.Range("O444").Value = 30000
$ws.Range("P444").Value = 30000
$ws.Range("Q444").Value = "`$/bandeja 10 kilos"
$ws.Range("R444").Value = "Perú"
$ws.Range("S444").Value = 3000
$ws.Range("T444").Value = 10

# Row 445
$ws.Range("D445").Value = 44162
$ws.Range("L445").Value = "Primera"
$ws.Range("M445").Value = 150
$ws.Range("N445").Value = 4200
$ws.Range("O445").Value = 4200
$ws.Range("P445").Value = 4200
$ws.Range("Q445").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R445").Value = "Provincia de Quillota"
$ws.Range("S445").Value = 4200
$ws.Range("T445").Value = 1

# Row 446
$ws.Range("D446").Value = 44162
$ws.Range("L446").Value = "Segunda"
$ws.Range("M446").Value = 150
$ws.Range("N446").Value = 3600
$ws.Range("O446").Value = 3600
$ws.Range("P446").Value = 3600
$ws.Range("Q446").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R446").Value = "Provincia de Quillota"
$ws.Range("S446").Value = 3600
$ws.Range("T446").Value = 1

# Row 447
$ws.Range("D447").Value = 44162
$ws.Range("L447").Value = "Tercera"
$ws.Range("M447").Value = 150
$ws.Range("N447").Value = 3000
$ws.Range("O447").Value = 3000
$ws.Range("P447").Value = 3000
$ws.Range("Q447").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R447").Value = "Provincia de Quillota"
$ws.Range("S447").Value = 3000
$ws.Range("T447").Value = 1

# Row 448
$ws.Range("D448").Value = 44235
$ws.Range("L448").Value = "Primera"
$ws.Range("M448").Value = 100
$ws.Range("N448").Value = 5300
$ws.Range("O448").Value = 5300
$ws.Range("P448").Value = 5300
$ws.Range("Q448").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R448").Value = "Provincia de Quillota"
$ws.Range("S448").Value = 5300
$ws.Range("T448").Value = 1

# Row 449
$ws.Range("D449").Value = 44726
$ws.Range("L449").Value = "Primera"
$ws.Range("M449").Value = 200
$ws.Range("N449").Value = 24000
$ws.Range("O449").Value = 25000
$ws.Range("P449").Value = 24500
$ws.Range("Q449").Value = "`$/bandeja 10 kilos"
$ws.Range("R449").Value = "Perú"
$ws.Range("S449").Value = 2450
$ws.Range("T449").Value = 10

# Row 450
$ws.Range("D450").Value = 44726
$ws.Range("L450").Value = "Primera"
$ws.Range("M450").Value = 300
$ws.Range("N450").Value = 4300
$ws.Range("O450").Value = 4400
$ws.Range("P450").Value = 4350
$ws.Range("Q450").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R450").Value = "Provincia de Quillota"
$ws.Range("S450").Value = 4350
$ws.Range("T450").Value = 1

# Row 451
$ws.Range("D451").Value = 44726
$ws.Range("L451").Value = "Segunda"
$ws.Range("M451").Value = 100
$ws.Range("N451").Value = 20000
$ws.Range("O451").Value = 20000
$ws.Range("P451").Value = 20000
$ws.Range("Q451").Value = "`$/bandeja 10 kilos"
$ws.Range("R451").Value = "Perú"
$ws.Range("S451").Value = 2000
$ws.Range("T451").Value = 10

# Row 452
$ws.Range("D452").Value = 44726
$ws.Range("L452").Value = "Segunda"
$ws.Range("M452").Value = 150
$ws.Range("N452").Value = 4000
$ws.Range("O452").Value = 4000
$ws.Range("P452").Value = 4000
$ws.Range("Q452").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R452").Value = "Provincia de Quillota"
$ws.Range("S452").Value = 4000
$ws.Range("T452").Value = 1

# Row 453
$ws.Range("D453").Value = 44662
$ws.Range("L453").Value = "Primera"
$ws.Range("M453").Value = 200
$ws.Range("N453").Value = 4300
$ws.Range("O453").Value = 4500
$ws.Range("P453").Value = 4400
$ws.Range("Q453").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R453").Value = "Provincia de Quillota"
$ws.Range("S453").Value = 4400
$ws.Range("T453").Value = 1

# Row 454
$ws.Range("D454").Value = 44662
$ws.Range("L454").Value = "Segunda"
$ws.Range("M454").Value = 100
$ws.Range("N454").Value = 3900
$ws.Range("O454").Value = 3900
$ws.Range("P454").Value = 3900
$ws.Range("Q454").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R454").Value = "Provincia de Quillota"
$ws.Range("S454").Value = 3900
$ws.Range("T454").Value = 1

# Row 455
$ws.Range("D455").Value = 44708
$ws.Range("L455").Value = "Primera"
$ws.Range("M455").Value = 250
$ws.Range("N455").Value = 4300
$ws.Range("O455").Value = 4300
$ws.Range("P455").Value = 4300
$ws.Range("Q455").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R455").Value = "Provincia de Quillota"
$ws.Range("S455").Value = 4300
$ws.Range("T455").Value = 1

# Row 456
$ws.Range("D456").Value = 44708
$ws.Range("L456").Value = "Segunda"
$ws.Range("M456").Value = 150
$ws.Range("N456").Value = 4000
$ws.Range("O456").Value = 4000
$ws.Range("P456").Value = 4000
$ws.Range("Q456").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R456").Value = "Provincia de Quillota"
$ws.Range("S456").Value = 4000
$ws.Range("T456").Value = 1

# Row 457
$ws.Range("D457").Value = 44628
$ws.Range("L457").Value = "Primera"
$ws.Range("M457").Value = 600
$ws.Range("N457").Value = 4300
$ws.Range("O457").Value = 4400
$ws.Range("P457").Value = 4350
$ws.Range("Q457").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R457").Value = "Provincia de Quillota"
$ws.Range("S457").Value = 4350
$ws.Range("T457").Value = 1

# Row 458
$ws.Range("D458").Value = 44483
$ws.Range("L458").Value = "1a nueva(o)"
$ws.Range("M458").Value = 200
$ws.Range("N458").Value = 4000
$ws.Range("O458").Value = 4200
$ws.Range("P458").Value = 4100
$ws.Range("Q458").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R458").Value = "Provincia de Quillota"
$ws.Range("S458").Value = 4100
$ws.Range("T458").Value = 1

# Row 459
$ws.Range("D459").Value = 44175
$ws.Range("L459").Value = "Primera"
$ws.Range("M459").Value = 100
$ws.Range("N459").Value = 4500
$ws.Range("O459").Value = 4500
$ws.Range("P459").Value = 4500
$ws.Range("Q459").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R459").Value = "Provincia de Quillota"
$ws.Range("S459").Value = 4500
$ws.Range("T459").Value = 1

# Row 460
$ws.Range("D460").Value = 44175
$ws.Range("L460").Value = "Segunda"
$ws.Range("M460").Value = 100
$ws.Range("N460").Value = 3800
$ws.Range("O460").Value = 3800
$ws.Range("P460").Value = 3800
$ws.Range("Q460").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R460").Value = "Provincia de Quillota"
$ws.Range("S460").Value = 3800
$ws.Range("T460").Value = 1

# Row 461
$ws.Range("D461").Value = 44175
$ws.Range("L461").Value = "Tercera"
$ws.Range("M461").Value = 60
$ws.Range("N461").Value = 3000
$ws.Range("O461").Value = 3000
$ws.Range("P461").Value = 3000
$ws.Range("Q461").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R461").Value = "Provincia de Quillota"
$ws.Range("S461").Value = 3000
$ws.Range("T461").Value = 1

# Row 462
$ws.Range("D462").Value = 44469
$ws.Range("L462").Value = "1a nueva(o)"
$ws.Range("M462").Value = 100
$ws.Range("N462").Value = 4000
$ws.Range("O462").Value = 4000
$ws.Range("P462").Value = 4000
$ws.Range("Q462").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R462").Value = "Provincia de Quillota"
$ws.Range("S462").Value = 4000
$ws.Range("T462").Value = 1

# Row 463
$ws.Range("D463").Value = 44469
$ws.Range("L463").Value = "2a nueva(o)"
$ws.Range("M463").Value = 100
$ws.Range("N463").Value = 3500
$ws.Range("O463").Value = 3500
$ws.Range("P463").Value = 3500
$ws.Range("Q463").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R463").Value = "Provincia de Quillota"
$ws.Range("S463").Value = 3500
$ws.Range("T463").Value = 1

# Row 464
$ws.Range("D464").Value = 44434
$ws.Range("L464").Value = "Especial"
$ws.Range("M464").Value = 100
$ws.Range("N464").Value = 35000
$ws.Range("O464").Value = 35000
$ws.Range("P464").Value = 35000
$ws.Range("Q464").Value = "`$/bandeja 10 kilos"
$ws.Range("R464").Value = "Perú"
$ws.Range("S464").Value = 3500
$ws.Range("T464").Value = 10

# Row 465
$ws.Range("D465").Value = 44434
$ws.Range("L465").Value = "Primera"
$ws.Range("M465").Value = 50
$ws.Range("N465").Value = 30000
$ws.Range("O465").Value = 30000
$ws.Range("P465").Value = 30000
$ws.Range("Q465").Value = "`$/bandeja 10 kilos"
$ws.Range("R465").Value = "Perú"
$ws.Range("S465").Value = 3000
$ws.Range("T465").Value = 10

# Row 466
$ws.Range("D466").Value = 44253
$ws.Range("L466").Value = "Primera"
$ws.Range("M466").Value = 200
$ws.Range("N466").Value = 5400
$ws.Range("O466").Value = 5500
$ws.Range("P466").Value = 5450
$ws.Range("Q466").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R466").Value = "Provincia de Quillota"
$ws.Range("S466").Value = 5450
$ws.Range("T466").Value = 1

# Row 467
$ws.Range("D467").Value = 44253
$ws.Range("L467").Value = "Segunda"
$ws.Range("M467").Value = 100
$ws.Range("N467").Value = 4300
$ws.Range("O467").Value = 4300
$ws.Range("P467").Value = 4300
$ws.Range("Q467").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R467").Value = "Provincia de Quillota"
$ws.Range("S467").Value = 4300
$ws.Range("T467").Value = 1

# Row 468
$ws.Range("D468").Value = 44204
$ws.Range("L468").Value = "Primera"
$ws.Range("M468").Value = 200
$ws.Range("N468").Value = 5400
$ws.Range("O468").Value = 5600
$ws.Range("P468").Value = 5500
$ws.Range("Q468").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R468").Value = "Provincia de Quillota"
$ws.Range("S468").Value = 5500
$ws.Range("T468").Value = 1

# Row 469
$ws.Range("D469").Value = 44204
$ws.Range("L469").Value = "Segunda"
$ws.Range("M469").Value = 100
$ws.Range("N469").Value = 4500
$ws.Range("O469").Value = 4500
$ws.Range("P469").Value = 4500
$ws.Range("Q469").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R469").Value = "Provincia de Quillota"
$ws.Range("S469").Value = 4500
$ws.Range("T469").Value = 1

# Row 470
$ws.Range("D470").Value = 44484
$ws.Range("L470").Value = "1a nueva(o)"
$ws.Range("M470").Value = 300
$ws.Range("N470").Value = 4000
$ws.Range("O470").Value = 4200
$ws.Range("P470").Value = 4100
$ws.Range("Q470").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R470").Value = "Provincia de Quillota"
$ws.Range("S470").Value = 4100
$ws.Range("T470").Value = 1

# Row 471
$ws.Range("D471").Value = 44484
$ws.Range("L471").Value = "2a nueva(o)"
$ws.Range("M471").Value = 150
$ws.Range("N471").Value = 3600
$ws.Range("O471").Value = 3600
$ws.Range("P471").Value = 3600
$ws.Range("Q471").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R471").Value = "Provincia de Quillota"
$ws.Range("S471").Value = 3600
$ws.Range("T471").Value = 1

# Row 472
$ws.Range("D472").Value = 44229
$ws.Range("L472").Value = "Primera"
$ws.Range("M472").Value = 200
$ws.Range("N472").Value = 5400
$ws.Range("O472").Value = 5500
$ws.Range("P472").Value = 5450
$ws.Range("Q472").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R472").Value = "Provincia de Quillota"
$ws.Range("S472").Value = 5450
$ws.Range("T472").Value = 1

# Row 473
$ws.Range("D473").Value = 44229
$ws.Range("L473").Value = "Segunda"
$ws.Range("M473").Value = 100
$ws.Range("N473").Value = 4200
$ws.Range("O473").Value = 4200
$ws.Range("P473").Value = 4200
$ws.Range("Q473").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R473").Value = "Provincia de Quillota"
$ws.Range("S473").Value = 4200
$ws.Range("T473").Value = 1

# Row 474
$ws.Range("D474").Value = 44231
$ws.Range("L474").Value = "Primera"
$ws.Range("M474").Value = 200
$ws.Range("N474").Value = 5400
$ws.Range("O474").Value = 5500
$ws.Range("P474").Value = 5450
$ws.Range("Q474").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R474").Value = "Provincia de Quillota"
$ws.Range("S474").Value = 5450
$ws.Range("T474").Value = 1

# Row 475
$ws.Range("D475").Value = 44231
$ws.Range("L475").Value = "Segunda"
$ws.Range("M475").Value = 100
$ws.Range("N475").Value = 4200
$ws.Range("O475").Value = 4200
$ws.Range("P475").Value = 4200
$ws.Range("Q475").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R475").Value = "Provincia de Quillota"
$ws.Range("S475").Value = 4200
$ws.Range("T475").Value = 1

# Row 476
$ws.Range("D476").Value = 44336
$ws.Range("L476").Value = "Primera"
$ws.Range("M476").Value = 120
$ws.Range("N476").Value = 7000
$ws.Range("O476").Value = 7200
$ws.Range("P476").Value = 7100
$ws.Range("Q476").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R476").Value = "Provincia de Quillota"
$ws.Range("S476").Value = 7100
$ws.Range("T476").Value = 1

# Row 477
$ws.Range("D477").Value = 44488
$ws.Range("L477").Value = "1a nueva(o)"
$ws.Range("M477").Value = 300
$ws.Range("N477").Value = 4000
$ws.Range("O477").Value = 4200
$ws.Range("P477").Value = 4100
$ws.Range("Q477").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R477").Value = "Provincia de Quillota"
$ws.Range("S477").Value = 4100
$ws.Range("T477").Value = 1

# Row 478
$ws.Range("D478").Value = 44488
$ws.Range("L478").Value = "2a nueva(o)"
$ws.Range("M478").Value = 150
$ws.Range("N478").Value = 3600
$ws.Range("O478").Value = 3600
$ws.Range("P478").Value = 3600
$ws.Range("Q478").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R478").Value = "Provincia de Quillota"
$ws.Range("S478").Value = 3600
$ws.Range("T478").Value = 1

# Row 479
$ws.Range("D479").Value = 44196
$ws.Range("L479").Value = "Primera"
$ws.Range("M479").Value = 200
$ws.Range("N479").Value = 5500
$ws.Range("O479").Value = 5600
$ws.Range("P479").Value = 5550
$ws.Range("Q479").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R479").Value = "Provincia de Quillota"
$ws.Range("S479").Value = 5550
$ws.Range("T479").Value = 1

# Row 480
$ws.Range("D480").Value = 44196
$ws.Range("L480").Value = "Segunda"
$ws.Range("M480").Value = 100
$ws.Range("N480").Value = 4600
$ws.Range("O480").Value = 4600
$ws.Range("P480").Value = 4600
$ws.Range("Q480").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R480").Value = "Provincia de Quillota"
$ws.Range("S480").Value = 4600
$ws.Range("T480").Value = 1

# Row 481
$ws.Range("D481").Value = 44369
$ws.Range("L481").Value = "1a nueva(o)"
$ws.Range("M481").Value = 300
$ws.Range("N481").Value = 5000
$ws.Range("O481").Value = 5200
$ws.Range("P481").Value = 5100
$ws.Range("Q481").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R481").Value = "Provincia de Quillota"
$ws.Range("S481").Value = 5100
$ws.Range("T481").Value = 1

# Row 482
$ws.Range("D482").Value = 44369
$ws.Range("L482").Value = "2a nueva(o)"
$ws.Range("M482").Value = 150
$ws.Range("N482").Value = 4200
$ws.Range("O482").Value = 4200
$ws.Range("P482").Value = 4200
$ws.Range("Q482").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R482").Value = "Provincia de Quillota"
$ws.Range("S482").Value = 4200
$ws.Range("T482").Value = 1

# Row 483
$ws.Range("D483").Value = 44369
$ws.Range("L483").Value = "Primera"
$ws.Range("M483").Value = 150
$ws.Range("N483").Value = 40000
$ws.Range("O483").Value = 40000
$ws.Range("P483").Value = 40000
$ws.Range("Q483").Value = "`$/bandeja 10 kilos"
$ws.Range("R483").Value = "Perú"
$ws.Range("S483").Value = 4000
$ws.Range("T483").Value = 10

# Row 484
$ws.Range("D484").Value = 44298
$ws.Range("L484").Value = "Primera"
$ws.Range("M484").Value = 80
$ws.Range("N484").Value = 6400
$ws.Range("O484").Value = 6500
$ws.Range("P484").Value = 6450
$ws.Range("Q484").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R484").Value = "Provincia de Quillota"
$ws.Range("S484").Value = 6450
$ws.Range("T484").Value = 1

# Row 485
$ws.Range("D485").Value = 44596
$ws.Range("L485").Value = "Primera"
$ws.Range("M485").Value = 400
$ws.Range("N485").Value = 4000
$ws.Range("O485").Value = 4100
$ws.Range("P485").Value = 4050
$ws.Range("Q485").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R485").Value = "Provincia de Quillota"
$ws.Range("S485").Value = 4050
$ws.Range("T485").Value = 1

# Row 486
$ws.Range("D486").Value = 44596
$ws.Range("L486").Value = "Segunda"
$ws.Range("M486").Value = 200
$ws.Range("N486").Value = 3500
$ws.Range("O486").Value = 3500
$ws.Range("P486").Value = 3500
$ws.Range("Q486").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R486").Value = "Provincia de Quillota"
$ws.Range("S486").Value = 3500
$ws.Range("T486").Value = 1

# Row 487
$ws.Range("D487").Value = 44399
$ws.Range("L487").Value = "Primera"
$ws.Range("M487").Value = 100
$ws.Range("N487").Value = 36000
$ws.Range("O487").Value = 36000
$ws.Range("P487").Value = 36000
$ws.Range("Q487").Value = "`$/bandeja 10 kilos"
$ws.Range("R487").Value = "Perú"
$ws.Range("S487").Value = 3600
$ws.Range("T487").Value = 10

# Row 488
$ws.Range("D488").Value = 44425
$ws.Range("L488").Value = "Especial"
$ws.Range("M488").Value = 100
$ws.Range("N488").Value = 35000
$ws.Range("O488").Value = 35000
$ws.Range("P488").Value = 35000
$ws.Range("Q488").Value = "`$/bandeja 10 kilos"
$ws.Range("R488").Value = "Perú"
$ws.Range("S488").Value = 3500
$ws.Range("T488").Value = 10

# Row 489
$ws.Range("D489").Value = 44425
$ws.Range("L489").Value = "Primera"
$ws.Range("M489").Value = 100
$ws.Range("N489").Value = 32000
$ws.Range("O489").Value = 32000
$ws.Range("P489").Value = 32000
$ws.Range("Q489").Value = "`$/bandeja 10 kilos"
$ws.Range("R489").Value = "Perú"
$ws.Range("S489").Value = 3200
$ws.Range("T489").Value = 10

# Row 490
$ws.Range("D490").Value = 44425
$ws.Range("L490").Value = "Segunda"
$ws.Range("M490").Value = 100
$ws.Range("N490").Value = 24000
$ws.Range("O490").Value = 24000
$ws.Range("P490").Value = 24000
$ws.Range("Q490").Value = "`$/bandeja 10 kilos"
$ws.Range("R490").Value = "Perú"
$ws.Range("S490").Value = 2400
$ws.Range("T490").Value = 10

# Row 491
$ws.Range("D491").Value = 44512
$ws.Range("L491").Value = "Primera"
$ws.Range("M491").Value = 400
$ws.Range("N491").Value = 4000
$ws.Range("O491").Value = 4200
$ws.Range("P491").Value = 4100
$ws.Range("Q491").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R491").Value = "Provincia de Quillota"
$ws.Range("S491").Value = 4100
$ws.Range("T491").Value = 1

# Row 492
$ws.Range("D492").Value = 44512
$ws.Range("L492").Value = "Segunda"
$ws.Range("M492").Value = 200
$ws.Range("N492").Value = 3600
$ws.Range("O492").Value = 3600
$ws.Range("P492").Value = 3600
$ws.Range("Q492").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R492").Value = "Provincia de Quillota"
$ws.Range("S492").Value = 3600
$ws.Range("T492").Value = 1

# Row 493
$ws.Range("D493").Value = 44397
$ws.Range("L493").Value = "Especial"
$ws.Range("M493").Value = 200
$ws.Range("N493").Value = 40000
$ws.Range("O493").Value = 40000
$ws.Range("P493").Value = 40000
$ws.Range("Q493").Value = "`$/bandeja 10 kilos"
$ws.Range("R493").Value = "Perú"
$ws.Range("S493").Value = 4000
$ws.Range("T493").Value = 10

# Row 494
$ws.Range("D494").Value = 44181
$ws.Range("L494").Value = "Primera"
$ws.Range("M494").Value = 40
$ws.Range("N494").Value = 4800
$ws.Range("O494").Value = 4800
$ws.Range("P494").Value = 4800
$ws.Range("Q494").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R494").Value = "Provincia de Quillota"
$ws.Range("S494").Value = 4800
$ws.Range("T494").Value = 1

# Row 495
$ws.Range("D495").Value = 44181
$ws.Range("L495").Value = "Segunda"
$ws.Range("M495").Value = 40
$ws.Range("N495").Value = 4000
$ws.Range("O495").Value = 4000
$ws.Range("P495").Value = 4000
$ws.Range("Q495").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R495").Value = "Provincia de Quillota"
$ws.Range("S495").Value = 4000
$ws.Range("T495").Value = 1

# Row 496
$ws.Range("D496").Value = 44497
$ws.Range("L496").Value = "1a nueva(o)"
$ws.Range("M496").Value = 200
$ws.Range("N496").Value = 3900
$ws.Range("O496").Value = 4000
$ws.Range("P496").Value = 3950
$ws.Range("Q496").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R496").Value = "Provincia de Quillota"
$ws.Range("S496").Value = 3950
$ws.Range("T496").Value = 1

# Row 497
$ws.Range("D497").Value = 44285
$ws.Range("L497").Value = "Primera"
$ws.Range("M497").Value = 200
$ws.Range("N497").Value = 6300
$ws.Range("O497").Value = 6400
$ws.Range("P497").Value = 6350
$ws.Range("Q497").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R497").Value = "Provincia de Quillota"
$ws.Range("S497").Value = 6350
$ws.Range("T497").Value = 1

# Row 498
$ws.Range("D498").Value = 44285
$ws.Range("L498").Value = "Segunda"
$ws.Range("M498").Value = 100
$ws.Range("N498").Value = 4990
$ws.Range("O498").Value = 4990
$ws.Range("P498").Value = 4990
$ws.Range("Q498").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R498").Value = "Provincia de Quillota"
$ws.Range("S498").Value = 4990
$ws.Range("T498").Value = 1

# Row 499
$ws.Range("D499").Value = 44362
$ws.Range("L499").Value = "1a nueva(o)"
$ws.Range("M499").Value = 150
$ws.Range("N499").Value = 5400
$ws.Range("O499").Value = 5400
$ws.Range("P499").Value = 5400
$ws.Range("Q499").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R499").Value = "Provincia de Quillota"
$ws.Range("S499").Value = 5400
$ws.Range("T499").Value = 1

# Row 500
$ws.Range("D500").Value = 44362
$ws.Range("L500").Value = "2a nueva(o)"
$ws.Range("M500").Value = 300
$ws.Range("N500").Value = 4500
$ws.Range("O500").Value = 5300
$ws.Range("P500").Value = 4900
$ws.Range("Q500").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R500").Value = "Provincia de Quillota"
$ws.Range("S500").Value = 4900
$ws.Range("T500").Value = 1

# Row 501
$ws.Range("D501").Value = 44362
$ws.Range("L501").Value = "Primera"
$ws.Range("M501").Value = 200
$ws.Range("N501").Value = 42000
$ws.Range("O501").Value = 43000
$ws.Range("P501").Value = 42500
$ws.Range("Q501").Value = "`$/bandeja 10 kilos"
$ws.Range("R501").Value = "Perú"
$ws.Range("S501").Value = 4250
$ws.Range("T501").Value = 10

# Row 502
$ws.Range("D502").Value = 44557
$ws.Range("L502").Value = "Primera"
$ws.Range("M502").Value = 400
$ws.Range("N502").Value = 3900
$ws.Range("O502").Value = 4000
$ws.Range("P502").Value = 3950
$ws.Range("Q502").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R502").Value = "Provincia de Quillota"
$ws.Range("S502").Value = 3950
$ws.Range("T502").Value = 1

# Row 503
$ws.Range("D503").Value = 44557
$ws.Range("L503").Value = "Segunda"
$ws.Range("M503").Value = 100
$ws.Range("N503").Value = 3500
$ws.Range("O503").Value = 3500
$ws.Range("P503").Value = 3500
$ws.Range("Q503").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R503").Value = "Provincia de Quillota"
$ws.Range("S503").Value = 3500
$ws.Range("T503").Value = 1

# Row 504
$ws.Range("D504").Value = 44747
$ws.Range("L504").Value = "Especial"
$ws.Range("M504").Value = 200
$ws.Range("N504").Value = 4800
$ws.Range("O504").Value = 4800
$ws.Range("P504").Value = 4800
$ws.Range("Q504").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R504").Value = "Provincia de Quillota"
$ws.Range("S504").Value = 4800
$ws.Range("T504").Value = 1

# Row 505
$ws.Range("D505").Value = 44747
$ws.Range("L505").Value = "Primera"
$ws.Range("M505").Value = 300
$ws.Range("N505").Value = 26000
$ws.Range("O505").Value = 26000
$ws.Range("P505").Value = 26000
$ws.Range("Q505").Value = "`$/bandeja 10 kilos"
$ws.Range("R505").Value = "Perú"
$ws.Range("S505").Value = 2600
$ws.Range("T505").Value = 10

# Row 506
$ws.Range("D506").Value = 44747
$ws.Range("L506").Value = "Primera"
$ws.Range("M506").Value = 150
$ws.Range("N506").Value = 4500
$ws.Range("O506").Value = 4500
$ws.Range("P506").Value = 4500
$ws.Range("Q506").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R506").Value = "Provincia de Quillota"
$ws.Range("S506").Value = 4500
$ws.Range("T506").Value = 1

# Row 507
$ws.Range("D507").Value = 44279
$ws.Range("L507").Value = "Primera"
$ws.Range("M507").Value = 80
$ws.Range("N507").Value = 6300
$ws.Range("O507").Value = 6400
$ws.Range("P507").Value = 6350
$ws.Range("Q507").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R507").Value = "Provincia de Quillota"
$ws.Range("S507").Value = 6350
$ws.Range("T507").Value = 1

# Row 508
$ws.Range("D508").Value = 44551
$ws.Range("L508").Value = "Primera"
$ws.Range("M508").Value = 400
$ws.Range("N508").Value = 3900
$ws.Range("O508").Value = 4000
$ws.Range("P508").Value = 3950
$ws.Range("Q508").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R508").Value = "Provincia de Quillota"
$ws.Range("S508").Value = 3950
$ws.Range("T508").Value = 1

# Row 509
$ws.Range("D509").Value = 44551
$ws.Range("L509").Value = "Segunda"
$ws.Range("M509").Value = 200
$ws.Range("N509").Value = 3500
$ws.Range("O509").Value = 3500
$ws.Range("P509").Value = 3500
$ws.Range("Q509").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R509").Value = "Provincia de Quillota"
$ws.Range("S509").Value = 3500
$ws.Range("T509").Value = 1

# Row 510
$ws.Range("D510").Value = 44517
$ws.Range("L510").Value = "Primera"
$ws.Range("M510").Value = 120
$ws.Range("N510").Value = 4000
$ws.Range("O510").Value = 4100
$ws.Range("P510").Value = 4050
$ws.Range("Q510").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R510").Value = "Provincia de Quillota"
$ws.Range("S510").Value = 4050
$ws.Range("T510").Value = 1

# Row 511
$ws.Range("D511").Value = 44517
$ws.Range("L511").Value = "Segunda"
$ws.Range("M511").Value = 40
$ws.Range("N511").Value = 3500
$ws.Range("O511").Value = 3500
$ws.Range("P511").Value = 3500
$ws.Range("Q511").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R511").Value = "Provincia de Quillota"
$ws.Range("S511").Value = 3500
$ws.Range("T511").Value = 1

# Row 512
$ws.Range("D512").Value = 44757
$ws.Range("L512").Value = "Especial"
$ws.Range("M512").Value = 70
$ws.Range("N512").Value = 6000
$ws.Range("O512").Value = 6000
$ws.Range("P512").Value = 6000
$ws.Range("Q512").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R512").Value = "Provincia de Quillota"
$ws.Range("S512").Value = 6000
$ws.Range("T512").Value = 1

# Row 513
$ws.Range("D513").Value = 44757
$ws.Range("L513").Value = "Primera"
$ws.Range("M513").Value = 500
$ws.Range("N513").Value = 28000
$ws.Range("O513").Value = 29000
$ws.Range("P513").Value = 28500
$ws.Range("Q513").Value = "`$/bandeja 10 kilos"
$ws.Range("R513").Value = "Perú"
$ws.Range("S513").Value = 2850
$ws.Range("T513").Value = 10

# Row 514
$ws.Range("D514").Value = 44757
$ws.Range("L514").Value = "Primera"
$ws.Range("M514").Value = 70
$ws.Range("N514").Value = 5500
$ws.Range("O514").Value = 5500
$ws.Range("P514").Value = 5500
$ws.Range("Q514").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R514").Value = "Provincia de Quillota"
$ws.Range("S514").Value = 5500
$ws.Range("T514").Value = 1

# Row 515
$ws.Range("D515").Value = 44547
$ws.Range("L515").Value = "Primera"
$ws.Range("M515").Value = 400
$ws.Range("N515").Value = 3900
$ws.Range("O515").Value = 4000
$ws.Range("P515").Value = 3950
$ws.Range("Q515").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R515").Value = "Provincia de Quillota"
$ws.Range("S515").Value = 3950
$ws.Range("T515").Value = 1

# Row 516
$ws.Range("D516").Value = 44547
$ws.Range("L516").Value = "Segunda"
$ws.Range("M516").Value = 200
$ws.Range("N516").Value = 3500
$ws.Range("O516").Value = 3500
$ws.Range("P516").Value = 3500
$ws.Range("Q516").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R516").Value = "Provincia de Quillota"
$ws.Range("S516").Value = 3500
$ws.Range("T516").Value = 1

# Row 517
$ws.Range("D517").Value = 44355
$ws.Range("L517").Value = "Primera"
$ws.Range("M517").Value = 200
$ws.Range("N517").Value = 42000
$ws.Range("O517").Value = 43000
$ws.Range("P517").Value = 42500
$ws.Range("Q517").Value = "`$/bandeja 10 kilos"
$ws.Range("R517").Value = "Perú"
$ws.Range("S517").Value = 4250
$ws.Range("T517").Value = 10

# Row 518
$ws.Range("D518").Value = 44657
$ws.Range("L518").Value = "Primera"
$ws.Range("M518").Value = 80
$ws.Range("N518").Value = 4300
$ws.Range("O518").Value = 4500
$ws.Range("P518").Value = 4400
$ws.Range("Q518").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R518").Value = "Provincia de Quillota"
$ws.Range("S518").Value = 4400
$ws.Range("T518").Value = 1

# Row 519
$ws.Range("D519").Value = 44186
$ws.Range("L519").Value = "Primera"
$ws.Range("M519").Value = 60
$ws.Range("N519").Value = 5000
$ws.Range("O519").Value = 5000
$ws.Range("P519").Value = 5000
$ws.Range("Q519").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R519").Value = "Provincia de Quillota"
$ws.Range("S519").Value = 5000
$ws.Range("T519").Value = 1

# Row 520
$ws.Range("D520").Value = 44186
$ws.Range("L520").Value = "Segunda"
$ws.Range("M520").Value = 60
$ws.Range("N520").Value = 4300
$ws.Range("O520").Value = 4300
$ws.Range("P520").Value = 4300
$ws.Range("Q520").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R520").Value = "Provincia de Quillota"
$ws.Range("S520").Value = 4300
$ws.Range("T520").Value = 1

# Row 521
$ws.Range("D521").Value = 44189
$ws.Range("L521").Value = "Primera"
$ws.Range("M521").Value = 400
$ws.Range("N521").Value = 5300
$ws.Range("O521").Value = 5400
$ws.Range("P521").Value = 5350
$ws.Range("Q521").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R521").Value = "Provincia de Quillota"
$ws.Range("S521").Value = 5350
$ws.Range("T521").Value = 1

# Row 522
$ws.Range("D522").Value = 44189
$ws.Range("L522").Value = "Tercera"
$ws.Range("M522").Value = 200
$ws.Range("N522").Value = 4500
$ws.Range("O522").Value = 4500
$ws.Range("P522").Value = 4500
$ws.Range("Q522").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R522").Value = "Provincia de Quillota"
$ws.Range("S522").Value = 4500
$ws.Range("T522").Value = 1

# Row 523
$ws.Range("D523").Value = 44609
$ws.Range("L523").Value = "Primera"
$ws.Range("M523").Value = 200
$ws.Range("N523").Value = 4000
$ws.Range("O523").Value = 4000
$ws.Range("P523").Value = 4000
$ws.Range("Q523").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R523").Value = "Provincia de Quillota"
$ws.Range("S523").Value = 4000
$ws.Range("T523").Value = 1

# Row 524
$ws.Range("D524").Value = 44609
$ws.Range("L524").Value = "Segunda"
$ws.Range("M524").Value = 100
$ws.Range("N524").Value = 3500
$ws.Range("O524").Value = 3500
$ws.Range("P524").Value = 3500
$ws.Range("Q524").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R524").Value = "Provincia de Quillota"
$ws.Range("S524").Value = 3500
$ws.Range("T524").Value = 1

# Row 525
$ws.Range("D525").Value = 44358
$ws.Range("L525").Value = "1a nueva(o)"
$ws.Range("M525").Value = 200
$ws.Range("N525").Value = 5300
$ws.Range("O525").Value = 5400
$ws.Range("P525").Value = 5350
$ws.Range("Q525").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R525").Value = "Provincia de Quillota"
$ws.Range("S525").Value = 5350
$ws.Range("T525").Value = 1

# Row 526
$ws.Range("D526").Value = 44358
$ws.Range("L526").Value = "2a nueva(o)"
$ws.Range("M526").Value = 100
$ws.Range("N526").Value = 4500
$ws.Range("O526").Value = 4500
$ws.Range("P526").Value = 4500
$ws.Range("Q526").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R526").Value = "Provincia de Quillota"
$ws.Range("S526").Value = 4500
$ws.Range("T526").Value = 1

# Row 527
$ws.Range("A527").Value = 4
$ws.Range("B527").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C527").Value = "Los Lagos"
$ws.Range("E527").Value = 10
$ws.Range("F527").Value = "Fruta"
$ws.Range("G527").Value = 100106
$ws.Range("H527").Value = "Oleaginosos"
$ws.Range("I527").Value = 100106002
$ws.Range("J527").Value = "Palta"
$ws.Range("K527").Value = "Hass"
$ws.Range("D527").Value = 44358
$ws.Range("L527").Value = "Primera"
$ws.Range("M527").Value = 100
$ws.Range("N527").Value = 42000
$ws.Range("O527").Value = 42000
$ws.Range("P527").Value = 42000
$ws.Range("Q527").Value = "`$/bandeja 10 kilos"
$ws.Range("R527").Value = "Perú"
$ws.Range("S527").Value = 4200
$ws.Range("T527").Value = 10

# Row 528
$ws.Range("A528").Value = 4
$ws.Range("B528").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C528").Value = "Los Lagos"
$ws.Range("E528").Value = 10
$ws.Range("F528").Value = "Fruta"
$ws.Range("G528").Value = 100106
$ws.Range("H528").Value = "Oleaginosos"
$ws.Range("I528").Value = 100106002
$ws.Range("J528").Value = "Palta"
$ws.Range("K528").Value = "Hass"
$ws.Range("D528").Value = 44572
$ws.Range("L528").Value = "Primera"
$ws.Range("M528").Value = 400
$ws.Range("N528").Value = 4000
$ws.Range("O528").Value = 4100
$ws.Range("P528").Value = 4050
$ws.Range("Q528").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R528").Value = "Provincia de Quillota"
$ws.Range("S528").Value = 4050
$ws.Range("T528").Value = 1

# Row 529
$ws.Range("A529").Value = 4
$ws.Range("B529").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C529").Value = "Los Lagos"
$ws.Range("E529").Value = 10
$ws.Range("F529").Value = "Fruta"
$ws.Range("G529").Value = 100106
$ws.Range("H529").Value = "Oleaginosos"
$ws.Range("I529").Value = 100106002
$ws.Range("J529").Value = "Palta"
$ws.Range("K529").Value = "Hass"
$ws.Range("D529").Value = 44572
$ws.Range("L529").Value = "Segunda"
$ws.Range("M529").Value = 200
$ws.Range("N529").Value = 3600
$ws.Range("O529").Value = 3600
$ws.Range("P529").Value = 3600
$ws.Range("Q529").Value = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R529").Value = "Provincia de Quillota"
$ws.Range("S529").Value = 3600
$ws.Range("T529").Value = 1

# Preserve date style for newly added rows in column D
$dateFormat = $ws.Range("D526").NumberFormat
$ws.Range("D527").NumberFormat = $dateFormat
$ws.Range("D528").NumberFormat = $dateFormat
$ws.Range("D529").NumberFormat = $dateFormat
